# "Fixes to GetTransactionData due to missing first record."
#
# The three "Exceptions e-mail" config rows (ExceptionsEmailSubject,
# ExceptionsEmailBody, ExceptionsEmailSender) were living on the wrong
# sheet - they need to move from Settings!A10:C12 down to
# Constants!A12:C14 (GetTransactionData was reading the first Constants
# record and silently dropping it because these extra rows weren't part
# of that table).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Settings")
$ws2 = $wb.Worksheets.Item("Constants")

# Move the data (values + formatting) from Settings to Constants.
$ws1.Range("A10:C12").Copy()
$ws2.Range("A12:C14").PasteSpecial()

# Remove the now-duplicated rows from their old location.
$ws1.Range("A10:C12").Clear()

# Reflect the new location in each sheet's selection, and make Constants
# (where the data now lives) the active tab.
$ws1.Range("A11:C12").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("A13:C14").Select() | Out-Null
